# Auto-generated Excel COM-interop script
# Applies updated currentAveragePrice / Leve price / profit values
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 76
$ws.Range("H76").Value = 9373130
$ws.Range("I76").Value = 3208
$ws.Range("K76").Value = 3208
$ws.Range("M76").Value = -2893
# Row 79
$ws.Range("H79").Value = 9373130
$ws.Range("I79").Value = 3208
$ws.Range("K79").Value = 3208
$ws.Range("M79").Value = -2116
# Row 80
$ws.Range("H80").Value = 1203757.6
$ws.Range("I80").Value = 2853562.5
$ws.Range("J80").Value = 3899.7273
$ws.Range("K80").Value = 8560687.5
$ws.Range("L80").Value = 11699.1819
$ws.Range("M80").Value = -8559689.5
$ws.Range("N80").Value = -13695.1819
# Row 83
$ws.Range("H83").Value = 1203757.6
$ws.Range("I83").Value = 2853562.5
$ws.Range("J83").Value = 3899.7273
$ws.Range("K83").Value = 25682062.5
$ws.Range("L83").Value = 35097.5457
$ws.Range("M83").Value = -25677070.5
$ws.Range("N83").Value = -45081.5457
# Row 101
$ws.Range("H101").Value = 414.25
$ws.Range("I101").Value = 435.66666
$ws.Range("K101").Value = 1306.99998
$ws.Range("M101").Value = 315.0000199999999
# Row 134
$ws.Range("H134").Value = 64748.5
$ws.Range("J134").Value = 64748.5
$ws.Range("L134").Value = 64748.5
$ws.Range("N134").Value = -74888.5
# Row 137
$ws.Range("H137").Value = 2503.6875
$ws.Range("I137").Value = 1774.625
$ws.Range("J137").Value = 2868.2188
$ws.Range("K137").Value = 5323.875
$ws.Range("L137").Value = 8604.6564
$ws.Range("M137").Value = -2773.875
$ws.Range("N137").Value = -13704.6564
# Row 138
$ws.Range("H138").Value = 3553.7
$ws.Range("I138").Value = 2435.4285
$ws.Range("J138").Value = 4155.846
$ws.Range("K138").Value = 7306.2855
$ws.Range("L138").Value = 12467.538
$ws.Range("M138").Value = -2166.2855
$ws.Range("N138").Value = -22747.538

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2779.6458
$ws.Range("I32").Value = 2150.1555
$ws.Range("J32").Value = 12222
$ws.Range("K32").Value = 2150.1555
$ws.Range("L32").Value = 12222
$ws.Range("M32").Value = -1863.1555
$ws.Range("N32").Value = -12796
# Row 37
$ws.Range("H37").Value = 11450
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents() | Out-Null
# Row 44
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents() | Out-Null
# Row 55
$ws.Range("H55").Value = 59997.5
$ws.Range("I55").Value = 59997.5
$ws.Range("K55").Value = 59997.5
$ws.Range("M55").Value = -59682.5
# Row 80
$ws.Range("H80").Value = 55073.934
$ws.Range("J80").Value = 59191.816
$ws.Range("L80").Value = 59191.816
$ws.Range("N80").Value = -61187.816
# Row 83
$ws.Range("H83").Value = 55073.934
$ws.Range("J83").Value = 59191.816
$ws.Range("L83").Value = 177575.448
$ws.Range("N83").Value = -187559.448

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 35571.285
$ws.Range("J82").Value = 45999.668
$ws.Range("L82").Value = 45999.668
$ws.Range("N82").Value = -46765.668
# Row 85
$ws.Range("H85").Value = 35571.285
$ws.Range("J85").Value = 45999.668
$ws.Range("L85").Value = 45999.668
$ws.Range("N85").Value = -48651.668
# Row 132
$ws.Range("H132").Value = 88664.836
$ws.Range("J132").Value = 88664.836
$ws.Range("L132").Value = 88664.836
$ws.Range("N132").Value = -98784.836
# Row 135
$ws.Range("H135").Value = 45217.832
$ws.Range("J135").Value = 45217.832
$ws.Range("L135").Value = 45217.832
$ws.Range("N135").Value = -55357.832

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1638.6296
$ws.Range("I58").Value = 1114.4736
$ws.Range("J58").Value = 2883.5
$ws.Range("K58").Value = 1114.4736
$ws.Range("L58").Value = 2883.5
$ws.Range("M58").Value = -911.4736
$ws.Range("N58").Value = -3289.5
# Row 135
$ws.Range("H135").Value = 56643.5
$ws.Range("I135").Value = 40000
$ws.Range("J135").Value = 59972.2
$ws.Range("K135").Value = 40000
$ws.Range("L135").Value = 59972.2
$ws.Range("M135").Value = -34930
$ws.Range("N135").Value = -70112.2
# Row 136
$ws.Range("H136").Value = 1638.6296
$ws.Range("I136").Value = 1114.4736
$ws.Range("J136").Value = 2883.5
$ws.Range("K136").Value = 3343.4208
$ws.Range("L136").Value = 8650.5
$ws.Range("M136").Value = -793.4207999999999
$ws.Range("N136").Value = -13750.5

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 6907.0586
$ws.Range("J131").Value = 11811.111
$ws.Range("L131").Value = 35433.333
$ws.Range("N131").Value = -45513.333
# Row 133
$ws.Range("H133").Value = 6302.0586
$ws.Range("I133").Value = 4903.8887
$ws.Range("K133").Value = 14711.6661
$ws.Range("M133").Value = -9651.666100000002
# Row 134
$ws.Range("H134").Value = 1980.9259
$ws.Range("I134").Value = 1353.625
$ws.Range("K134").Value = 4060.875
$ws.Range("M134").Value = 1009.125

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 4515.8887
$ws.Range("I43").Value = 4515.8887
$ws.Range("K43").Value = 4515.8887
$ws.Range("M43").Value = -4364.8887
# Row 46
$ws.Range("H46").Value = 3001
$ws.Range("I46").Value = 3001
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 3001
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -2845
$ws.Range("N46").ClearContents() | Out-Null
# Row 57
$ws.Range("H57").Value = 25008.125
$ws.Range("I57").Value = 14295
$ws.Range("J57").Value = 100000
$ws.Range("K57").Value = 14295
$ws.Range("L57").Value = 100000
$ws.Range("M57").Value = -13475
$ws.Range("N57").Value = -101640
# Row 133
$ws.Range("H133").Value = 60999.4
$ws.Range("J133").Value = 60999.4
$ws.Range("L133").Value = 60999.4
$ws.Range("N133").Value = -71119.39999999999

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 5559.75
$ws.Range("I68").Value = 5746.3335
$ws.Range("J68").Value = 5000
$ws.Range("K68").Value = 5746.3335
$ws.Range("L68").Value = 5000
$ws.Range("M68").Value = -4997.3335
$ws.Range("N68").Value = -6498
# Row 71
$ws.Range("H71").Value = 5559.75
$ws.Range("I71").Value = 5746.3335
$ws.Range("J71").Value = 5000
$ws.Range("K71").Value = 28731.6675
$ws.Range("L71").Value = 25000
$ws.Range("M71").Value = -24987.6675
$ws.Range("N71").Value = -32488
# Row 122
$ws.Range("H122").Value = 4337.154
$ws.Range("I122").Value = 3336.6333
$ws.Range("J122").Value = 7672.222
$ws.Range("K122").Value = 10009.8999
$ws.Range("L122").Value = 23016.666
$ws.Range("M122").Value = -7559.8999
$ws.Range("N122").Value = -27916.666
# Row 135
$ws.Range("H135").Value = 46249.75
$ws.Range("J135").Value = 46249.75
$ws.Range("L135").Value = 46249.75
$ws.Range("N135").Value = -56389.75

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 135
$ws.Range("H135").Value = 54666.168
$ws.Range("J135").Value = 54666.168
$ws.Range("L135").Value = 54666.168
$ws.Range("N135").Value = -64806.168

